$wb = $excel.ActiveWorkbook
$headerSrc = $wb.Worksheets.Item("newcastle_upon_tyne_properties")

# ---- Devon_properties ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Devon_properties"
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "average_rating"
$ws.Range("C1").Value = "number_of_ratings"
$ws.Range("D1").Value = "price_per_night"
# Reuse the bold/bordered/centered header style already in the workbook
# (style index 1, shared by every other sheet's header row) instead of
# defining a new one.
$headerSrc.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$ws.Range("A2").Value = "No Snakes on thi..."
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 147
$ws.Range("A3").Value = "Exceptionally be..."
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 53
$ws.Range("D3").Value = 88
$ws.Range("A4").Value = "Tranquil Room by..."
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 59
$ws.Range("A5").Value = "Cosy barn betwee..."
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 82
$ws.Range("A6").Value = "The Cabin Devon ..."
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 409
$ws.Range("D6").Value = 122
$ws.Range("A7").Value = "The Cabin at Axe..."
$ws.Range("B7").Value = 4.91
$ws.Range("C7").Value = 150
$ws.Range("D7").Value = 114
$ws.Range("A8").Value = "Hattie - luxury ..."
$ws.Range("B8").Value = 4.97
$ws.Range("C8").Value = 72
$ws.Range("D8").Value = 71
$ws.Range("A9").Value = "Swaledale Shephe..."
$ws.Range("B9").Value = 4.82
$ws.Range("C9").Value = 198
$ws.Range("D9").Value = 59
$ws.Range("A10").Value = "Cosy Idyllic Cab..."
$ws.Range("B10").Value = 4.87
$ws.Range("C10").Value = 551
$ws.Range("D10").Value = 78
$ws.Range("A11").Value = "Lakeside Lodge, ..."
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 119
$ws.Range("D11").Value = 155
$ws.Range("A12").Value = "The Wizards Caul..."
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 165
$ws.Range("D12").Value = 176
$ws.Range("A13").Value = "Luxury Roundhous..."
$ws.Range("B13").Value = 4.97
$ws.Range("C13").Value = 38
$ws.Range("D13").Value = 206
$ws.Range("A14").Value = "Outstanding self..."
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 381
$ws.Range("D14").Value = 117
$ws.Range("A15").Value = "Cosy traditional..."
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 64
$ws.Range("A16").Value = "Luxury Shepherds..."
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 130
$ws.Range("D16").Value = 159
$ws.Range("A17").Value = "The Shippon. Uni..."
$ws.Range("B17").Value = 4.96
$ws.Range("C17").Value = 331
$ws.Range("D17").Value = 212
$ws.Range("A18").Value = "ModburyLittleHom..."
$ws.Range("B18").Value = 4.84
$ws.Range("C18").Value = 133
$ws.Range("D18").Value = 76

# ---- Bristol_properties ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Bristol_properties"
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "average_rating"
$ws.Range("C1").Value = "number_of_ratings"
$ws.Range("D1").Value = "price_per_night"
# Reuse the bold/bordered/centered header style already in the workbook
# (style index 1, shared by every other sheet's header row) instead of
# defining a new one.
$headerSrc.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$ws.Range("A2").Value = "Modern Comfort C..."
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 42
$ws.Range("A3").Value = "The Great Room N..."
$ws.Range("B3").Value = 4.91
$ws.Range("C3").Value = 94
$ws.Range("D3").Value = 45
$ws.Range("A4").Value = "Toad Lodge The B..."
$ws.Range("B4").Value = 4.83
$ws.Range("C4").Value = 345
$ws.Range("D4").Value = 58
$ws.Range("A5").Value = "Comfortable, cos..."
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 55
$ws.Range("A6").Value = "Entire home in E..."
$ws.Range("D6").Value = 76
$ws.Range("A7").Value = "Twin room in fam..."
$ws.Range("B7").Value = 4.88
$ws.Range("C7").Value = 130
$ws.Range("D7").Value = 47
$ws.Range("A8").Value = "Clifton cosy bed..."
$ws.Range("B8").Value = 4.78
$ws.Range("C8").Value = 181
$ws.Range("D8").Value = 33
$ws.Range("A9").Value = "Comfortable room..."
$ws.Range("B9").Value = 4.98
$ws.Range("C9").Value = 166
$ws.Range("D9").Value = 39
$ws.Range("A10").Value = "Spacious Double ..."
$ws.Range("B10").Value = 4.95
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = 42
$ws.Range("A11").Value = "Cosy boutique ci..."
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = 100
$ws.Range("A12").Value = "Perfectly Locate..."
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 117
$ws.Range("A13").Value = "Tuscany House..."
$ws.Range("B13").Value = 4.91
$ws.Range("C13").Value = 45
$ws.Range("D13").Value = 60
$ws.Range("A14").Value = "Palm House an ur..."
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 29
$ws.Range("D14").Value = 43
$ws.Range("A15").Value = "The Robinson Bui..."
$ws.Range("B15").Value = 4.93
$ws.Range("C15").Value = 135
$ws.Range("D15").Value = 96
$ws.Range("A16").Value = "Private self con..."
$ws.Range("B16").Value = 4.92
$ws.Range("C16").Value = 145
$ws.Range("D16").Value = 71
$ws.Range("A17").Value = "Riverside Walk..."
$ws.Range("B17").Value = 4.85
$ws.Range("C17").Value = 213
$ws.Range("D17").Value = 88
$ws.Range("A18").Value = "Self contained a..."
$ws.Range("B18").Value = 4.96
$ws.Range("C18").Value = 141
$ws.Range("D18").Value = 75

# ---- birmingham_properties ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "birmingham_properties"
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "average_rating"
$ws.Range("C1").Value = "number_of_ratings"
$ws.Range("D1").Value = "price_per_night"
# Reuse the bold/bordered/centered header style already in the workbook
# (style index 1, shared by every other sheet's header row) instead of
# defining a new one.
$headerSrc.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$ws.Range("A2").Value = "Plough House - 5..."
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 210
$ws.Range("D2").Value = 179
$ws.Range("A3").Value = "Double room with..."
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 30
$ws.Range("A4").Value = "1 Bed Flat, with..."
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 99
$ws.Range("A5").Value = "Comfy NEC/Airpor..."
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 102
$ws.Range("A6").Value = "Kula Birmingham ..."
$ws.Range("B6").Value = 4.81
$ws.Range("C6").Value = 113
$ws.Range("D6").Value = 131
$ws.Range("A7").Value = "Double Room2 wit..."
$ws.Range("B7").Value = 4.92
$ws.Range("C7").Value = 253
$ws.Range("D7").Value = 28
$ws.Range("A8").Value = "Cozy Stay Near A..."
$ws.Range("B8").Value = 4.83
$ws.Range("C8").Value = 46
$ws.Range("D8").Value = 33
$ws.Range("A9").Value = "The Foxes Den - ..."
$ws.Range("B9").Value = 4.96
$ws.Range("C9").Value = 620
$ws.Range("D9").Value = 69
$ws.Range("A10").Value = "''Heron's Rest' c..."
$ws.Range("B10").Value = 4.92
$ws.Range("C10").Value = 162
$ws.Range("D10").Value = 108
$ws.Range("A11").Value = "Lovely Room by U..."
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 30
$ws.Range("A12").Value = "Unique chalet ho..."
$ws.Range("B12").Value = 4.96
$ws.Range("C12").Value = 363
$ws.Range("D12").Value = 81
$ws.Range("A13").Value = "2 Bed Flat - Cen..."
$ws.Range("D13").Value = 80
$ws.Range("A14").Value = "Cosy 1 bedroom e..."
$ws.Range("B14").Value = 4.96
$ws.Range("C14").Value = 28
$ws.Range("D14").Value = 48
$ws.Range("A15").Value = "The Blue Room..."
$ws.Range("B15").Value = 4.88
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 35
$ws.Range("A16").Value = "Studio Near HS2,..."
$ws.Range("D16").Value = 71
$ws.Range("A17").Value = "Room in West Mid..."
$ws.Range("B17").Value = 4.93
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 31
$ws.Range("A18").Value = "Bed 10m from Bir..."
$ws.Range("B18").Value = 4.9
$ws.Range("C18").Value = 81
$ws.Range("D18").Value = 32

# ---- new_york_properties ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "new_york_properties"
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "average_rating"
$ws.Range("C1").Value = "number_of_ratings"
$ws.Range("D1").Value = "price_per_night"
# Reuse the bold/bordered/centered header style already in the workbook
# (style index 1, shared by every other sheet's header row) instead of
# defining a new one.
$headerSrc.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$ws.Range("A2").Value = "Home from home &..."
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 113
$ws.Range("A3").Value = "3 BD | Luxury St..."
$ws.Range("D3").Value = 95
$ws.Range("A4").Value = "Holywell Grange ..."
$ws.Range("B4").Value = 4.92
$ws.Range("C4").Value = 162
$ws.Range("D4").Value = 69
$ws.Range("A5").Value = "Home from home,b..."
$ws.Range("B5").Value = 4.86
$ws.Range("C5").Value = 166
$ws.Range("D5").Value = 51
$ws.Range("A6").Value = "3 bed house in S..."
$ws.Range("D6").Value = 136
$ws.Range("A7").Value = "A single room 20..."
$ws.Range("B7").Value = 4.91
$ws.Range("C7").Value = 108
$ws.Range("D7").Value = 39
$ws.Range("A8").Value = "1 Bedroom House ..."
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 135
$ws.Range("D8").Value = 123
$ws.Range("A9").Value = "Stylish 3 bed ho..."
$ws.Range("B9").Value = 4.88
$ws.Range("C9").Value = 77
$ws.Range("D9").Value = 133
$ws.Range("A10").Value = "Comfortable sing..."
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 37
$ws.Range("A11").Value = "Seghill's Sanctu..."
$ws.Range("B11").Value = 4.89
$ws.Range("C11").Value = 123
$ws.Range("D11").Value = 76
$ws.Range("A12").Value = "RestfullStays- S..."
$ws.Range("B12").Value = 4.64
$ws.Range("C12").Value = 33
$ws.Range("D12").Value = 109
$ws.Range("A13").Value = "Stylish and Cosy..."
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 84
$ws.Range("A14").Value = "Tiny Homestead@W..."
$ws.Range("B14").Value = 4.99
$ws.Range("C14").Value = 212
$ws.Range("D14").Value = 126
$ws.Range("A15").Value = "House in Westmoo..."
$ws.Range("B15").Value = 4.98
$ws.Range("C15").Value = 42
$ws.Range("D15").Value = 133
$ws.Range("A16").Value = "THE PLUMES Heato..."
$ws.Range("B16").Value = 4.96
$ws.Range("C16").Value = 285
$ws.Range("D16").Value = 61
$ws.Range("A17").Value = "Studio in leafy ..."
$ws.Range("B17").Value = 4.98
$ws.Range("C17").Value = 485
$ws.Range("D17").Value = 67
$ws.Range("A18").Value = "Lovely bright si..."
$ws.Range("B18").Value = 4.99
$ws.Range("C18").Value = 167
$ws.Range("D18").Value = 38

